# Use conventional styles/indents for Word bullet lists.
#
# 1) Every abstractNum level's hanging indent goes from 480 -> 360
#    twips (the conventional Word default hanging indent).
# 2) The bullet list's (abstractNumId 991) per-level glyphs switch
#    from the plain "bullet/en-dash" glyphs to the conventional
#    Word "3 bullets" gallery: Symbol bullet, Courier New "o",
#    Wingdings square, cycling across the 9 levels, each carrying an
#    explicit <w:rFonts .../> hint so the glyph renders with the
#    right symbol font.

$d = $word.ActiveDocument
$full = $d.Content.WordOpenXML

# Private-use-area glyphs used by the Symbol / Wingdings fonts for
# the conventional Word bullet gallery.
$symbolBullet   = [char]0xF0B7
$wingdingSquare = [char]0xF0A7

# --- Step 1: global hanging indent 480 -> 360 (all abstractNums) ---
$full = $full.Replace('w:hanging="480"', 'w:hanging="360"')

# --- Step 2: rewrite the 9 levels of the bullet abstractNum (991) ---
# Build old/new level fragments and replace them one at a time so we
# don't disturb the (already-edited) 480->360 text elsewhere.

$levels = @(
    @{ Ilvl = 0; Left = 720;  OldGlyph = [char]0x2022; NewGlyph = $symbolBullet;   Font = 'Symbol' },
    @{ Ilvl = 1; Left = 1440; OldGlyph = [char]0x2013; NewGlyph = 'o';             Font = 'Courier New' },
    @{ Ilvl = 2; Left = 2160; OldGlyph = [char]0x2022; NewGlyph = $wingdingSquare; Font = 'Wingdings' },
    @{ Ilvl = 3; Left = 2880; OldGlyph = [char]0x2013; NewGlyph = $symbolBullet;   Font = 'Symbol' },
    @{ Ilvl = 4; Left = 3600; OldGlyph = [char]0x2022; NewGlyph = 'o';             Font = 'Courier New' },
    @{ Ilvl = 5; Left = 4320; OldGlyph = [char]0x2013; NewGlyph = $wingdingSquare; Font = 'Wingdings' },
    @{ Ilvl = 6; Left = 5040; OldGlyph = [char]0x2022; NewGlyph = $symbolBullet;   Font = 'Symbol' },
    @{ Ilvl = 7; Left = 5760; OldGlyph = [char]0x2013; NewGlyph = 'o';             Font = 'Courier New' },
    @{ Ilvl = 8; Left = 6480; OldGlyph = [char]0x2022; NewGlyph = $wingdingSquare; Font = 'Wingdings' }
)

foreach ($lvl in $levels) {
    $old = '<w:lvl w:ilvl="' + $lvl.Ilvl + '"><w:numFmt w:val="bullet" /><w:lvlText w:val="' + $lvl.OldGlyph + '" /><w:lvlJc w:val="left" /><w:pPr><w:ind w:left="' + $lvl.Left + '" w:hanging="360" /></w:pPr></w:lvl>'

    $new = '<w:lvl w:ilvl="' + $lvl.Ilvl + '"><w:numFmt w:val="bullet" /><w:lvlText w:val="' + $lvl.NewGlyph + '" /><w:lvlJc w:val="left" /><w:pPr><w:ind w:left="' + $lvl.Left + '" w:hanging="360" /><w:rFonts w:ascii="' + $lvl.Font + '" w:hAnsi="' + $lvl.Font + '" w:cs="' + $lvl.Font + '" w:hint="default" /></w:pPr></w:lvl>'

    if ($full.IndexOf($old) -lt 0) {
        throw "Could not locate expected level fragment for ilvl=$($lvl.Ilvl)"
    }

    $full = $full.Replace($old, $new)
}

$d.Content.WordOpenXML = $full

Write-Output "done"
